$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("C1").Value = "rules"
$ws.Range("E1").Value = "adaptive_filter"

# Column E becomes a text label "RLS" for all data rows (was numeric 1)
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 5).Value = "RLS"
}

# Updated (recomputed) metric values for F (RMSE), G (NDEI), H (MAE)
$values = @{
    2  = @(0.706743346563334, 0.6449055915693407, 0.5879232502272569)
    3  = @(0.6772034483415389, 0.6179503388170352, 0.4994116277870033)
    4  = @(0.6433792482569667, 0.5870856467460416, 0.4419407946614844)
    5  = @(0.07297556577197656, 0.06659043999937597, 0.06416105955444051)
    6  = @(0.04652094794186699, 0.04245051559198183, 0.03632379677429636)
    7  = @(0.04609559993322777, 0.04206238415718561, 0.03583260967931559)
    8  = @(0.04576221002271383, 0.04175816478461034, 0.03534115031537344)
    9  = @(0.04549875227985394, 0.04151775874139992, 0.03489693720560434)
    10 = @(0.04528951787868621, 0.04132683167301238, 0.03448736321811946)
    11 = @(0.04512342619218241, 0.04117527247141489, 0.03407073794079594)
    12 = @(0.04499288425335137, 0.04105615252964845, 0.03370464516948388)
    13 = @(0.04489304331708718, 0.04096504735210784, 0.03338685953756324)
    14 = @(0.04482133512792257, 0.04089961339736536, 0.03306429629208981)
    15 = @(0.04477719633850795, 0.04085933660914089, 0.03278703734952997)
    16 = @(0.04476192591403198, 0.040845402297374, 0.0325065917460539)
}

foreach ($r in $values.Keys) {
    $row = $values[$r]
    $ws.Cells.Item($r, 6).Value = $row[0]
    $ws.Cells.Item($r, 7).Value = $row[1]
    $ws.Cells.Item($r, 8).Value = $row[2]
}
